# Practice Quizzes are uploaded.
# Adds a computed "Question Concatenated" column to the Astronomy sheet
# (Count + Question, via CONCAT), tweaks some view/row/column cosmetics
# on the other sheets, and leaves Remote_Sensing as the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Astronomy: insert a new column C = CONCAT(A, ". ", B) ("Question
# Concatenated"), pushing the old Answer column from C to D.
# ---------------------------------------------------------------------
$astro = $wb.Worksheets.Item("Astronomy")

$astro.Columns.Item(3).Insert()

$astro.Range("C1").Value = "Question Concatenated"
$astro.Range("C2:C21").Formula = "=_xlfn.CONCAT(A2,"". "",B2)"

# Column widths: B:D all become a uniform 35.33203125, replacing the old
# bestFit B (69.88671875) / custom C (42.77734375) widths.
$astro.Range("B1:D1").ColumnWidth = 35.33203125

# Row heights (auto-fit result of the narrower wrapped columns).
$astro.Rows.Item(2).RowHeight = 28.8
$astro.Rows.Item(3).RowHeight = 43.2
$astro.Rows.Item(4).RowHeight = 28.8
$astro.Rows.Item(5).RowHeight = 57.6
$astro.Rows.Item(6).RowHeight = 43.2
$astro.Rows.Item(7).RowHeight = 28.8
$astro.Rows.Item(8).RowHeight = 43.2
$astro.Rows.Item(9).RowHeight = 28.8
$astro.Rows.Item(10).RowHeight = 28.8
$astro.Rows.Item(11).RowHeight = 43.2
$astro.Rows.Item(12).RowHeight = 28.8
$astro.Rows.Item(13).RowHeight = 72
$astro.Range("A14:A21").EntireRow.RowHeight = 28.8

# View: unfreeze the prior "tabSelected" state, scroll the frozen pane
# down, and move the active cell to D13.
$astro.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$astro.Range("D13").Select()

# ---------------------------------------------------------------------
# Dynamic_Planet: cosmetic-only changes -- wider B:C columns, a couple
# of row heights, and a different scroll/selection position.
# ---------------------------------------------------------------------
$dynPlanet = $wb.Worksheets.Item("Dynamic_Planet")

$dynPlanet.Range("B1:C1").ColumnWidth = 53
$dynPlanet.Rows.Item(2).RowHeight = 28.8
$dynPlanet.Rows.Item(3).RowHeight = 115.2

$dynPlanet.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$dynPlanet.Range("B10").Select()

# ---------------------------------------------------------------------
# Remote_Sensing: becomes the active tab; scroll/selection change.
# ---------------------------------------------------------------------
$remote = $wb.Worksheets.Item("Remote_Sensing")

$remote.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$remote.Range("B12").Select()
